$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Ativação:" date-like text value (B8/C8) and the mirrored row 15 (B15/C15) ---
# These cells hold a text string that looks like a date ("01/01/2012" -> "01/01/2023").
# Force text format first so Excel doesn't auto-convert the string to a date serial number,
# then restore the original plain (non-numFmt) cell style by copying formats from a
# neighbouring plain-styled cell so the style index matches the rest of the column.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2023"
$ws.Range("C15").Value = "01/01/2023"
$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# --- Add new "Objectives:" text in row 11 (B11/C11), matching row 10's plain style ---
$objectivesText = "To present experimental techniques for the characterization of electrical, magnetic and thermal  properties of materials."
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("B11").Value = $objectivesText
$ws.Range("C11").Value = $objectivesText

# --- Same text is mirrored into row 14 (Short syllabus row), matching row 10's plain style ---
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("B14").Value = $objectivesText
$ws.Range("C14").Value = $objectivesText

# --- Add new "Syllabus:" text in row 16 (B16/C16), matching row 13's plain style ---
$syllabusText = "histerese de materiais magnéticos macios. Medidas de magnetostricção.Propriedades térmicas dos materiais:  expansão térmica.Electrical properties: electrical conductivity in pure metals, metallic alloys and semiconductors, and superconductors; Hall Effect; Ohm's Law and dependence on temperature.Magnetic properties: magnetic susceptibility and c.c. magnetization. Hysteresis curves of soft magnetic materials. Magnetostriction measurements.Thermal properties of materials: thermal expansion."
$ws.Range("B13").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("B16").Value = $syllabusText
$ws.Range("C16").Value = $syllabusText

# --- Update "Norma de recuperação:" value (B20/C20) ---
$normaText = "Média aritmética das notas dos relatórios de cada experimento"
$ws.Range("B20").Value = $normaText
$ws.Range("C20").Value = $normaText
